# Generate Report for Archive
#
# Refresh the localization-status report:
#   - the Status value "Ready for handoff" becomes "In Translation" for
#     the zh-cn and de-de rows (the Overview sheet mirrors each locale's
#     status in its own column, E for zh-cn and F for de-de).
#   - with the new, shorter status text in place, the Status-holding
#     columns are narrowed to fit the content.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newColumnWidth = 12.5

# --- Overview sheet: zh-cn status is column E, de-de status is column F ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth

# --- zh-cn sheet: Status is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth

# --- de-de sheet: Status is column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
